# "Moved Brine-belly to correct tab"
# Fills in the real TIER values for the Umbar Baharbêl deed rows (previously
# placeholder 0s) and adds the "Brine-belly" deed as a new data row, which
# pushes the trailing helper-formula row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Umbar Baharbêl")

# --- 1. Fill in the TIER (column C) values that were previously placeholder 0s ---
$tierUpdates = @{
    22 = 1;  23 = 2;  24 = 3;  25 = 4;  26 = 3;  27 = 4;  28 = 3;  29 = 4;
    30 = 2;  31 = 3;  32 = 4;  33 = 3;  34 = 4;  35 = 3;  36 = 4;
    37 = 2;  38 = 3;  39 = 4;  40 = 3;  41 = 4;  42 = 3;  43 = 4;
    44 = 2;  45 = 3;  46 = 4;  47 = 3;  48 = 4;  49 = 3;  50 = 4;
    51 = 2;  52 = 3;  53 = 4;  54 = 3;  55 = 4;
    56 = 1;  57 = 2;  58 = 2;  59 = 2;  60 = 2;  61 = 2;
    62 = 1;  63 = 2;  64 = 2;  65 = 2;  66 = 2;  67 = 2;
    69 = 1;  70 = 2;  71 = 1;  72 = 2;  73 = 1;  74 = 2;  75 = 1;  76 = 2;
}

foreach ($row in $tierUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $tierUpdates[$row]
}

# --- 2. Add the new "Brine-belly" deed row (ID 1879485102), moved in from another tab ---
$dataRow = 78
$ws.Cells.Item($dataRow, 1).Value = 1879485102
$ws.Cells.Item($dataRow, 2).Value = "Brine-belly"
$ws.Cells.Item($dataRow, 3).Value = 0

# --- 3. Re-build the helper-column formulas (G:N) for the data row and append
#        a fresh trailing blank-helper row beneath it, matching the sheet's
#        existing per-row formula pattern. ---
function Set-HelperFormulas($row) {
    $ws.Cells.Item($row, 7).Formula  = '=IF(AND(A' + $row + '>0,C' + $row + '>0),CONCATENATE("[",A' + $row + ',"] = ",C' + $row + ',", // ",B' + $row + '),"")'
    $ws.Cells.Item($row, 8).Formula  = '=CONCATENATE(J' + $row + ',K' + $row + ',L' + $row + ',N' + $row + '," -- ",B' + $row + ')'
    $ws.Cells.Item($row, 9).Formula  = '=ROW()-1'
    $ws.Cells.Item($row, 10).Formula = '=CONCATENATE(REPT(" ",2-LEN(I' + $row + ')),"[",I' + $row + ',"] = {")'
    $ws.Cells.Item($row, 11).Formula = '=IF(LEN(A' + $row + ')>0,CONCATENATE("[""ID""] = ",A' + $row + ',"; "),"")'
    $ws.Cells.Item($row, 12).Formula = '=IF(LEN(D' + $row + ')>0,CONCATENATE("[""CAT_ID""] = ",D' + $row + ',"; "),"")'
    $ws.Cells.Item($row, 13).Formula = '=CONCATENATE("[""TIER""] = ",TEXT(C' + $row + ',"0"),"; ")'
    $ws.Cells.Item($row, 14).Formula = '=CONCATENATE("};")'
}

Set-HelperFormulas 78
Set-HelperFormulas 79

# --- 4. Window / view bookkeeping ---
$excel.Width = 22770
$excel.Height = 15600
$excel.Left = 4575
$excel.Top = 0

$ws.Activate()
$ws.Range("D2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H9").Select()
